$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038545788421697
$ws.Range("D2").Value = 1.046034768426814
$ws.Range("E2").Value = 1.051740884680747
$ws.Range("F2").Value = 1.058391620907857
$ws.Range("I2").Value = 1.035594522201494
$ws.Range("J2").Value = 1.043642641997537
$ws.Range("K2").Value = 1.048801334733194
$ws.Range("L2").Value = 1.05449154759217
$ws.Range("M2").Value = 1.06112398211486
$ws.Range("N2").Value = 1.018445592372575

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039503037091454
$ws.Range("D3").Value = 1.046780739190166
$ws.Range("E3").Value = 1.052655513223995
$ws.Range("F3").Value = 1.059320223615432
$ws.Range("I3").Value = 1.035738884518591
$ws.Range("J3").Value = 1.044244689193413
$ws.Range("K3").Value = 1.049359037196958
$ws.Range("L3").Value = 1.0552186208268
$ws.Range("M3").Value = 1.061866318203179
$ws.Range("N3").Value = 1.018648242489969

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04012304373128
$ws.Range("D4").Value = 1.047263736957685
$ws.Range("E4").Value = 1.053248259563544
$ws.Range("F4").Value = 1.05992195193933
$ws.Range("I4").Value = 1.035831036957044
$ws.Range("J4").Value = 1.044634238688564
$ws.Range("K4").Value = 1.049719546619597
$ws.Range("L4").Value = 1.055689373726623
$ws.Range("M4").Value = 1.062346888634839
$ws.Range("N4").Value = 1.018779289798479

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040383837801498
$ws.Range("D5").Value = 1.047466860728518
$ws.Range("E5").Value = 1.053497669094569
$ws.Range("F5").Value = 1.060175123015449
$ws.Range("I5").Value = 1.035869475727683
$ws.Range("J5").Value = 1.044798000630667
$ws.Range("K5").Value = 1.049871017154528
$ws.Range("L5").Value = 1.055887345983951
$ws.Range("M5").Value = 1.062548974088989
$ws.Range("N5").Value = 1.018834362324698

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040427634645877
$ws.Range("D6").Value = 1.047500970254295
$ws.Range("E6").Value = 1.05353955885204
$ws.Range("F6").Value = 1.06021764350194
$ws.Range("I6").Value = 1.035875912034145
$ws.Range("J6").Value = 1.044825496678601
$ws.Range("K6").Value = 1.049896444552992
$ws.Range("L6").Value = 1.05592059033485
$ws.Range("M6").Value = 1.062582908235086
$ws.Range("N6").Value = 1.018843608073619

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040126527912687
$ws.Range("D7").Value = 1.047266450828305
$ws.Range("E7").Value = 1.05325159132726
$ws.Range("F7").Value = 1.059925334021973
$ws.Range("I7").Value = 1.035831551765575
$ws.Range("J7").Value = 1.044636426903238
$ws.Range("K7").Value = 1.049721570922008
$ws.Range("L7").Value = 1.055692018776935
$ws.Range("M7").Value = 1.062349588701362
$ws.Range("N7").Value = 1.018780025758426

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038869169676123
$ws.Range("D8").Value = 1.046286808749497
$ws.Range("E8").Value = 1.052049796237225
$ws.Range("F8").Value = 1.058705267533934
$ws.Range("I8").Value = 1.035643570603682
$ws.Range("J8").Value = 1.043846109230273
$ws.Range("K8").Value = 1.048989887044181
$ws.Range("L8").Value = 1.054737205071513
$ws.Range("M8").Value = 1.061374809826018
$ws.Range("N8").Value = 1.018514095431039

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036658205123553
$ws.Range("D9").Value = 1.044562953149357
$ws.Range("E9").Value = 1.049939191733037
$ws.Range("F9").Value = 1.056562007857986
$ws.Range("I9").Value = 1.035302701041656
$ws.Range("J9").Value = 1.042453395649674
$ws.Range("K9").Value = 1.047697848185785
$ws.Range("L9").Value = 1.053056958646867
$ws.Range("M9").Value = 1.059658937087342
$ws.Range("N9").Value = 1.018044890612857

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035187422298132
$ws.Range("D10").Value = 1.043415415389303
$ws.Range("E10").Value = 1.048536986590825
$ws.Range("F10").Value = 1.055137727866639
$ws.Range("I10").Value = 1.035069017753928
$ws.Range("J10").Value = 1.041524930056349
$ws.Range("K10").Value = 1.046834727106839
$ws.Range("L10").Value = 1.051938381983725
$ws.Range("M10").Value = 1.058516314760539
$ws.Range("N10").Value = 1.017731706649128

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034551326732329
$ws.Range("D11").Value = 1.042918939978701
$ws.Range("E11").Value = 1.047930987322392
$ws.Range("F11").Value = 1.054522098429376
$ws.Range("I11").Value = 1.034966310399446
$ws.Range("J11").Value = 1.04112290886999
$ws.Range("K11").Value = 1.046460581971679
$ws.Range("L11").Value = 1.051454417033012
$ws.Range("M11").Value = 1.058021868802369
$ws.Range("N11").Value = 1.017596009165598

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034315168096726
$ws.Range("D12").Value = 1.042734590841702
$ws.Range("E12").Value = 1.047706068577781
$ws.Range("F12").Value = 1.054293591912868
$ws.Range("I12").Value = 1.034927932333582
$ws.Range("J12").Value = 1.040973582844978
$ws.Range("K12").Value = 1.046321547585982
$ws.Range("L12").Value = 1.05127471005646
$ws.Range("M12").Value = 1.057838258411189
$ws.Range("N12").Value = 1.017545592431876

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034365819680475
$ws.Range("D13").Value = 1.0427741314411
$ws.Range("E13").Value = 1.047754306409131
$ws.Range("F13").Value = 1.05434259981662
$ws.Range("I13").Value = 1.034936174867197
$ws.Range("J13").Value = 1.041005613653713
$ws.Range("K13").Value = 1.046351373641648
$ws.Range("L13").Value = 1.051313255117708
$ws.Range("M13").Value = 1.057877641233081
$ws.Range("N13").Value = 1.017556407558444

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034531803422639
$ws.Range("D14").Value = 1.042903700308956
$ws.Range("E14").Value = 1.047912391863084
$ws.Range("F14").Value = 1.054503206621618
$ws.Range("I14").Value = 1.034963142704264
$ws.Range("J14").Value = 1.041110565476736
$ws.Range("K14").Value = 1.046449090568307
$ws.Range("L14").Value = 1.051439561184914
$ws.Range("M14").Value = 1.058006690506514
$ws.Range("N14").Value = 1.017591841958398

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03463408682775
$ws.Range("D15").Value = 1.042983540479817
$ws.Range("E15").Value = 1.048009816939815
$ws.Range("F15").Value = 1.054602183770112
$ws.Range("I15").Value = 1.034979728284688
$ws.Range("J15").Value = 1.041175230124533
$ws.Range("K15").Value = 1.046509289243851
$ws.Range("L15").Value = 1.051517390394351
$ws.Range("M15").Value = 1.058086208533032
$ws.Range("N15").Value = 1.017613672597539

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035229654023043
$ws.Range("D16").Value = 1.043448373714596
$ws.Range("E16").Value = 1.048577229424265
$ws.Range("F16").Value = 1.055178608319037
$ws.Range("I16").Value = 1.035075802102839
$ws.Range("J16").Value = 1.041551611182464
$ws.Range("K16").Value = 1.04685954938268
$ws.Range("L16").Value = 1.051970509376414
$ws.Range("M16").Value = 1.058549136294318
$ws.Range("N16").Value = 1.017740710649407

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035603442550679
$ws.Range("D17").Value = 1.043740063713463
$ws.Range("E17").Value = 1.048933465101881
$ws.Range("F17").Value = 1.055540477896514
$ws.Range("I17").Value = 1.035135659839043
$ws.Range("J17").Value = 1.041787708575813
$ws.Range("K17").Value = 1.047079149836583
$ws.Range("L17").Value = 1.052254843126363
$ws.Range("M17").Value = 1.058839604392416
$ws.Range("N17").Value = 1.017820375308966

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035821540640448
$ws.Range("D18").Value = 1.043910241597725
$ws.Range("E18").Value = 1.049141363641307
$ws.Range("F18").Value = 1.055751655680677
$ws.Range("I18").Value = 1.035170427059135
$ws.Range("J18").Value = 1.041925421143653
$ws.Range("K18").Value = 1.047207199626552
$ws.Range("L18").Value = 1.052420727370932
$ws.Range("M18").Value = 1.0590090600285
$ws.Range("N18").Value = 1.017866833959889

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035895918902626
$ws.Range("D19").Value = 1.043968274606782
$ws.Range("E19").Value = 1.049212270676434
$ws.Range("F19").Value = 1.055823679678979
$ws.Range("I19").Value = 1.035182256869564
$ws.Range("J19").Value = 1.041972377702905
$ws.Range("K19").Value = 1.047250854573508
$ws.Range("L19").Value = 1.052477295882892
$ws.Range("M19").Value = 1.059066845159012
$ws.Range("N19").Value = 1.017882673723131

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035563330963162
$ws.Range("D20").Value = 1.043708763984903
$ws.Range("E20").Value = 1.048895232759896
$ws.Range("F20").Value = 1.055501641805073
$ws.Range("I20").Value = 1.035129252849279
$ws.Range("J20").Value = 1.041762377454412
$ws.Range("K20").Value = 1.047055592857779
$ws.Range("L20").Value = 1.052224332952048
$ws.Range("M20").Value = 1.058808436749955
$ws.Range("N20").Value = 1.017811828910682

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03448292214953
$ws.Range("D21").Value = 1.042865543730321
$ws.Range("E21").Value = 1.047865834757534
$ws.Range("F21").Value = 1.054455907335869
$ws.Range("I21").Value = 1.034955207637591
$ws.Range("J21").Value = 1.041079659702781
$ws.Range("K21").Value = 1.046420317022763
$ws.Range("L21").Value = 1.051402365554067
$ws.Range("M21").Value = 1.057968687353574
$ws.Range("N21").Value = 1.017581407756775

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033804295421076
$ws.Range("D22").Value = 1.042335748313315
$ws.Range("E22").Value = 1.047219632452334
$ws.Range("F22").Value = 1.053799371777887
$ws.Range("I22").Value = 1.034844459735689
$ws.Range("J22").Value = 1.040650423148039
$ws.Range("K22").Value = 1.046020545988657
$ws.Range("L22").Value = 1.050885904845221
$ws.Range("M22").Value = 1.057440986831182
$ws.Range("N22").Value = 1.017436459855976

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034163984547213
$ws.Range("D23").Value = 1.042616567278001
$ws.Range("E23").Value = 1.047562099278132
$ws.Range("F23").Value = 1.054147322300582
$ws.Range("I23").Value = 1.034903294122689
$ws.Range("J23").Value = 1.040877967806171
$ws.Range("K23").Value = 1.04623250479554
$ws.Range("L23").Value = 1.051159657561259
$ws.Range("M23").Value = 1.05772070360576
$ws.Range("N23").Value = 1.017513306268282

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035581455435727
$ws.Range("D24").Value = 1.043722906861382
$ws.Range("E24").Value = 1.048912507962706
$ws.Range("F24").Value = 1.055519189837642
$ws.Range("I24").Value = 1.035132148345757
$ws.Range("J24").Value = 1.041773823494949
$ws.Range("K24").Value = 1.047066237364464
$ws.Range("L24").Value = 1.052238119071756
$ws.Range("M24").Value = 1.058822519971835
$ws.Range("N24").Value = 1.017815690686176

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037229233951417
$ws.Range("D25").Value = 1.045008317953137
$ws.Range("E25").Value = 1.050483982251491
$ws.Range("F25").Value = 1.057115294497856
$ws.Range("I25").Value = 1.035391960955286
$ws.Range("J25").Value = 1.0428134480417
$ws.Range("K25").Value = 1.048032186560517
$ws.Range("L25").Value = 1.053491067918914
$ws.Range("M25").Value = 1.060102308902541
$ws.Range("N25").Value = 1.018166260152694

Write-Output "Updated vm_pu values for rows 2-25 (380 kV case)"
